$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The daily job-matches sheet is refreshed with this run's results.
# Existing rows 2-18 are overwritten in place and rows 19-23 are newly
# appended, matching the latest "job_matches_2026-02-23" export.
# Columns: A=Title, B=Company, C=Location, D=Match Score (%),
#          E=Matched Keywords, F=Posted At, G=Apply Link

# Posted At (column F) holds plain date-stamp text (e.g. "2026-02-23").
# Force the whole column range to Text format first so Excel stores the
# values as literal strings instead of auto-converting them to date serials.
$postedRange = $ws.Range("F2:F23")
$postedRange.NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = 'GenAI Architect'
$ws.Cells.Item(2, 2).Value = 'Tata Consultancy Services (TCS)'
$ws.Cells.Item(2, 3).Value = 'Edison, NJ, US USA'
$ws.Cells.Item(2, 4).Value = 24.4
$ws.Cells.Item(2, 5).Value = 'AI Engineer, LangChain, RAG, LLaMA, Hugging Face, Pinecone, TensorFlow, PyTorch, S3, Data Lake'
$ws.Cells.Item(2, 6).Value = '2026-02-23'
$ws.Cells.Item(2, 7).Value = 'https://www.indeed.com/viewjob?jk=e014fe941f1f22ec'

$ws.Cells.Item(3, 1).Value = 'Senior Data Engineer (1043) - DataSF'
$ws.Cells.Item(3, 2).Value = 'City and County of San Francisco'
$ws.Cells.Item(3, 3).Value = 'San Francisco, CA, US USA'
$ws.Cells.Item(3, 4).Value = 16.7
$ws.Cells.Item(3, 5).Value = 'Data Scientist, RAG, BigQuery, Kinesis, Terraform, Snowflake, Databricks, BigQuery, PySpark, Kafka'
$ws.Cells.Item(3, 6).Value = '2026-02-23'
$ws.Cells.Item(3, 7).Value = 'https://www.indeed.com/viewjob?jk=012bacb5fe848397'

$ws.Cells.Item(4, 1).Value = 'DevOps Engineer'
$ws.Cells.Item(4, 2).Value = 'kp reddy'
$ws.Cells.Item(4, 3).Value = 'San Francisco, CA, US USA'
$ws.Cells.Item(4, 4).Value = 15.6
$ws.Cells.Item(4, 5).Value = 'RAG, S3, EC2, Docker, Kubernetes, CI/CD, Jenkins, GitHub Actions, Terraform, Git'
$ws.Cells.Item(4, 6).Value = '2026-02-23'
$ws.Cells.Item(4, 7).Value = 'https://www.indeed.com/viewjob?jk=31ba26cff0eff16a'

$ws.Cells.Item(5, 1).Value = 'DevOps Engineer'
$ws.Cells.Item(5, 2).Value = 'kp reddy'
$ws.Cells.Item(5, 3).Value = 'Atlanta, GA, US USA'
$ws.Cells.Item(5, 4).Value = 15.6
$ws.Cells.Item(5, 5).Value = 'RAG, S3, EC2, Docker, Kubernetes, CI/CD, Jenkins, GitHub Actions, Terraform, Git'
$ws.Cells.Item(5, 6).Value = '2026-02-23'
$ws.Cells.Item(5, 7).Value = 'https://www.indeed.com/viewjob?jk=dce2730f22d939ec'

$ws.Cells.Item(6, 1).Value = 'Senior Software Engineer New'
$ws.Cells.Item(6, 2).Value = 'Convey'
$ws.Cells.Item(6, 3).Value = 'Chicago, IL, US USA'
$ws.Cells.Item(6, 4).Value = 14.4
$ws.Cells.Item(6, 5).Value = 'RAG, Copilot, Kinesis, Docker, Kubernetes, Git, Kafka, MongoDB, NoSQL, SQL'
$ws.Cells.Item(6, 6).Value = '2026-02-23'
$ws.Cells.Item(6, 7).Value = 'https://www.indeed.com/viewjob?jk=f5324f4c6c59141c'

$ws.Cells.Item(7, 1).Value = 'Software Engineer New'
$ws.Cells.Item(7, 2).Value = 'Convey'
$ws.Cells.Item(7, 3).Value = 'Chicago, IL, US USA'
$ws.Cells.Item(7, 4).Value = 14.4
$ws.Cells.Item(7, 5).Value = 'RAG, Copilot, Kinesis, Docker, Kubernetes, Git, Kafka, MongoDB, NoSQL, SQL'
$ws.Cells.Item(7, 6).Value = '2026-02-23'
$ws.Cells.Item(7, 7).Value = 'https://www.indeed.com/viewjob?jk=121b212372d04928'

$ws.Cells.Item(8, 1).Value = 'Senior Software Engineer - AI, Building Design'
$ws.Cells.Item(8, 2).Value = 'kp reddy'
$ws.Cells.Item(8, 3).Value = 'Atlanta, GA, US USA'
$ws.Cells.Item(8, 4).Value = 14.4
$ws.Cells.Item(8, 5).Value = 'AI Engineer, Generative AI, PyTorch, YOLO, AWS SageMaker, Azure ML, MLflow, Docker, Kubernetes, Python'
$ws.Cells.Item(8, 6).Value = '2026-02-23'
$ws.Cells.Item(8, 7).Value = 'https://www.indeed.com/viewjob?jk=c1a9a570fd3e2cff'

$ws.Cells.Item(9, 1).Value = 'Senior Software Engineer - AI, Building Design'
$ws.Cells.Item(9, 2).Value = 'kp reddy'
$ws.Cells.Item(9, 3).Value = 'San Francisco, CA, US USA'
$ws.Cells.Item(9, 4).Value = 14.4
$ws.Cells.Item(9, 5).Value = 'AI Engineer, Generative AI, PyTorch, YOLO, AWS SageMaker, Azure ML, MLflow, Docker, Kubernetes, Python'
$ws.Cells.Item(9, 6).Value = '2026-02-23'
$ws.Cells.Item(9, 7).Value = 'https://www.indeed.com/viewjob?jk=30cf3c854ab272eb'

$ws.Cells.Item(10, 1).Value = 'Senior Software Engineer - AI, Building Design'
$ws.Cells.Item(10, 2).Value = 'kp reddy'
$ws.Cells.Item(10, 3).Value = 'San Francisco, CA, US USA'
$ws.Cells.Item(10, 4).Value = 14.4
$ws.Cells.Item(10, 5).Value = 'AI Engineer, Generative AI, PyTorch, YOLO, AWS SageMaker, Azure ML, MLflow, Docker, Kubernetes, Python'
$ws.Cells.Item(10, 6).Value = '2026-02-23'
$ws.Cells.Item(10, 7).Value = 'https://www.indeed.com/viewjob?jk=e50a761ac02df3ba'

$ws.Cells.Item(11, 1).Value = 'Software Engineer – CRG (Analyst / Associate)'
$ws.Cells.Item(11, 2).Value = 'Goldman Sachs'
$ws.Cells.Item(11, 3).Value = 'Dallas, TX, US USA'
$ws.Cells.Item(11, 4).Value = 14.4
$ws.Cells.Item(11, 5).Value = 'LangChain, RAG, Prompt Engineering, Kubernetes, CI/CD, Git, Kafka, MongoDB, Python, SQL'
$ws.Cells.Item(11, 6).Value = '2026-02-23'
$ws.Cells.Item(11, 7).Value = 'https://www.indeed.com/viewjob?jk=ca332a980655e235'

$ws.Cells.Item(12, 1).Value = 'Software Engineer - Full Stack'
$ws.Cells.Item(12, 2).Value = 'kp reddy'
$ws.Cells.Item(12, 3).Value = 'Atlanta, GA, US USA'
$ws.Cells.Item(12, 4).Value = 13.3
$ws.Cells.Item(12, 5).Value = 'RAG, Docker, CI/CD, Jenkins, GitHub Actions, Git, PostgreSQL, MySQL, Python, SQL'
$ws.Cells.Item(12, 6).Value = '2026-02-23'
$ws.Cells.Item(12, 7).Value = 'https://www.indeed.com/viewjob?jk=15451134c42bdf64'

$ws.Cells.Item(13, 1).Value = 'Software Engineer - Full Stack'
$ws.Cells.Item(13, 2).Value = 'kp reddy'
$ws.Cells.Item(13, 3).Value = 'San Francisco, CA, US USA'
$ws.Cells.Item(13, 4).Value = 13.3
$ws.Cells.Item(13, 5).Value = 'RAG, Docker, CI/CD, Jenkins, GitHub Actions, Git, PostgreSQL, MySQL, Python, SQL'
$ws.Cells.Item(13, 6).Value = '2026-02-23'
$ws.Cells.Item(13, 7).Value = 'https://www.indeed.com/viewjob?jk=a3288093ad8485f4'

$ws.Cells.Item(14, 1).Value = 'Data Engineer - ITS4'
$ws.Cells.Item(14, 2).Value = 'State of Minnesota - Minnesota IT Services'
$ws.Cells.Item(14, 3).Value = 'Saint Paul, MN, US USA'
$ws.Cells.Item(14, 4).Value = 13.3
$ws.Cells.Item(14, 5).Value = 'RAG, S3, Glue, Redshift, Data Lake, Git, Redshift, PySpark, NoSQL, Python'
$ws.Cells.Item(14, 6).Value = '2026-02-18'
$ws.Cells.Item(14, 7).Value = 'https://www.indeed.com/viewjob?jk=54c8deff754c16fa'

$ws.Cells.Item(15, 1).Value = 'Associate Data Scientist'
$ws.Cells.Item(15, 2).Value = 'MetLife'
$ws.Cells.Item(15, 3).Value = 'Cary, NC, US USA'
$ws.Cells.Item(15, 4).Value = 13.3
$ws.Cells.Item(15, 5).Value = 'Data Scientist, Generative AI, RAG, Copilot, Prompt Engineering, TensorFlow, PyTorch, spaCy, Git, Python'
$ws.Cells.Item(15, 6).Value = '2026-02-23'
$ws.Cells.Item(15, 7).Value = 'https://www.indeed.com/viewjob?jk=916958a0f37fc7ab'

$ws.Cells.Item(16, 1).Value = 'AI Platform Engineer'
$ws.Cells.Item(16, 2).Value = 'OrderlyMeds'
$ws.Cells.Item(16, 3).Value = 'Remote, US USA'
$ws.Cells.Item(16, 4).Value = 12.2
$ws.Cells.Item(16, 5).Value = 'AI Engineer, RAG, Kubernetes, CI/CD, GitHub Actions, Terraform, Git, Python, R, Scala'
$ws.Cells.Item(16, 6).Value = '2026-02-23'
$ws.Cells.Item(16, 7).Value = 'https://www.indeed.com/viewjob?jk=4bdae7064fd5dc6f'

$ws.Cells.Item(17, 1).Value = 'Applied AI Engineer'
$ws.Cells.Item(17, 2).Value = 'Propio LS LLC'
$ws.Cells.Item(17, 3).Value = 'Overland Park, KS, US USA'
$ws.Cells.Item(17, 4).Value = 12.2
$ws.Cells.Item(17, 5).Value = 'AI Engineer, LangChain, Hugging Face, FAISS, Pinecone, Prompt Engineering, FastAPI, Python, R, Scala'
$ws.Cells.Item(17, 6).Value = '2026-02-23'
$ws.Cells.Item(17, 7).Value = 'https://www.indeed.com/viewjob?jk=3055ad914d70d64e'

$ws.Cells.Item(18, 1).Value = 'AI Developer'
$ws.Cells.Item(18, 2).Value = 'Bourns'
$ws.Cells.Item(18, 3).Value = 'Carrollton, TX, US USA'
$ws.Cells.Item(18, 4).Value = 11.1
$ws.Cells.Item(18, 5).Value = 'LangChain, RAG, LLaMA, Prompt Engineering, FastAPI, Kubernetes, Git, Python, R, Scala'
$ws.Cells.Item(18, 6).Value = '2026-02-23'
$ws.Cells.Item(18, 7).Value = 'https://www.indeed.com/viewjob?jk=6ed224a6ba5ee451'

$ws.Cells.Item(19, 1).Value = 'Data Scientist Subcontractor'
$ws.Cells.Item(19, 2).Value = 'The Hackett Group'
$ws.Cells.Item(19, 3).Value = 'FL, US USA'
$ws.Cells.Item(19, 4).Value = 11.1
$ws.Cells.Item(19, 5).Value = 'Data Scientist, Generative AI, LangChain, LLaMA, Hugging Face, TensorFlow, PyTorch, XGBoost, Python, R'
$ws.Cells.Item(19, 6).Value = '2026-02-23'
$ws.Cells.Item(19, 7).Value = 'https://www.indeed.com/viewjob?jk=84b99efc2d74d07c'

$ws.Cells.Item(20, 1).Value = 'Jr. Software Developer'
$ws.Cells.Item(20, 2).Value = 'Buyers Edge Platform'
$ws.Cells.Item(20, 3).Value = 'Remote, US USA'
$ws.Cells.Item(20, 4).Value = 10
$ws.Cells.Item(20, 5).Value = 'RAG, S3, EC2, CI/CD, Git, MySQL, SQL, R, Java'
$ws.Cells.Item(20, 6).Value = '2026-02-23'
$ws.Cells.Item(20, 7).Value = 'https://www.indeed.com/viewjob?jk=90037e7aa244c1c1'

$ws.Cells.Item(21, 1).Value = 'Machine Learning Researcher'
$ws.Cells.Item(21, 2).Value = 'Rivet Industries'
$ws.Cells.Item(21, 3).Value = 'San Jose, CA, US USA'
$ws.Cells.Item(21, 4).Value = 10
$ws.Cells.Item(21, 5).Value = 'TensorFlow, PyTorch, Docker, Kubernetes, CI/CD, Python, R, Scala, Optimization'
$ws.Cells.Item(21, 6).Value = '2026-02-23'
$ws.Cells.Item(21, 7).Value = 'https://www.indeed.com/viewjob?jk=756a548f7b63644b'

$ws.Cells.Item(22, 1).Value = 'Compliance, Dallas, Associate, Software Engineering'
$ws.Cells.Item(22, 2).Value = 'Goldman Sachs'
$ws.Cells.Item(22, 3).Value = 'Dallas, TX, US USA'
$ws.Cells.Item(22, 4).Value = 10
$ws.Cells.Item(22, 5).Value = 'Data Scientist, RAG, MongoDB, NoSQL, Python, SQL, R, Java, Scala'
$ws.Cells.Item(22, 6).Value = '2026-02-23'
$ws.Cells.Item(22, 7).Value = 'https://www.indeed.com/viewjob?jk=ead7d1da7f34ba0a'

$ws.Cells.Item(23, 1).Value = 'MLOps Engineer'
$ws.Cells.Item(23, 2).Value = 'ValueBase Consulting'
$ws.Cells.Item(23, 3).Value = 'Ann Arbor, MI, US USA'
$ws.Cells.Item(23, 4).Value = 10
$ws.Cells.Item(23, 5).Value = 'Azure ML, Docker, Kubernetes, CI/CD, Jenkins, Python, R, Java, Scala'
$ws.Cells.Item(23, 6).Value = '2026-02-23'
$ws.Cells.Item(23, 7).Value = 'https://www.indeed.com/viewjob?jk=a39342bfa0110d52'

# Restore the default (style-less) cell formatting on the Posted At column
# now that the text values are locked in, so no stray number format lingers
# on the cells themselves.
$postedRange.Style = "Normal"

